$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously described the "Cases" tab; this test suite now covers
# Participants instead, so update the tab-name label accordingly.
$ws.Range("A2").Value = "ParticipantsTab"

# Reflect the cell the author left selected after making the edit.
$ws.Range("A2").Select()
